# Enigma workbook: add a "reflector" configuration sheet showing the
# rotor/reflector wiring (A<->B, C<->D, ... pairing) so it can eventually
# be changed.

$wb = $excel.ActiveWorkbook

# --- update the selection left on Sheet3 (now that its data is done) -----
$ws3 = $wb.Worksheets.Item(3)
[void]$ws3.Activate()
[void]$ws3.Range("I1:I26").Select()

# --- add the new "reflector" sheet at the very end of the workbook -------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "reflector"

# --- letters A..Z go down column B ---------------------------------------
$letters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")
for ($i = 0; $i -lt 26; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $letters[$i]
    $ws.Cells.Item($row, 3).Value = $i
}

# --- column D: reflector wiring (adjacent pairs swap) ---------------------
$ws.Range("D1").Value = 1
$ws.Range("D2").Formula = "=IF(ISODD(C2),C1,C3)"
$ws.Range("D3:D26").Formula = "=IF(ISODD(C3),C2,C4)"

# --- column E: running, comma-separated trace of column D -----------------
$ws.Range("E1").Formula = "=D1"
$ws.Range("E2").Formula = '=E1&", "&D2'
$ws.Range("E3").Formula = '=E2&", "&D3'
$ws.Range("E4:E26").Formula = '=E3&", "&D4'

# --- the new sheet becomes the active tab, with D2 selected ---------------
[void]$ws.Range("D2").Select()

$wb.Save()
